# Add the 2022-Q3 sheet (copy the existing 2022-Q2 sheet as a formatting
# template, then overwrite it with the 2022-Q3 figures) and update the
# "总计" (summary) sheet with a new leading row for 2022-Q3.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet right after "2022-Q2" (i.e. right
#    after "总计"), by copying "2022-Q2" so it inherits identical
#    formatting (header style, column-A style, borders, etc.), then
#    relabel it and replace its data.
# ---------------------------------------------------------------------
$templateWs = $wb.Worksheets.Item("2022-Q2")
$templateWs.Copy($templateWs, $null)
$newWs = $wb.Worksheets.Item(2)
$newWs.Name = "2022-Q3"

# Drop the extra rows copied from 2022-Q2 (it had 23 data rows, 2022-Q3
# only has 15), keeping the header row + formatting intact.
$newWs.Range("A17:H24").Delete()

$q3 = @(
    @("003713", "英大睿盛灵活配置混合A", "2.83", "93.65", "7.09", "0.2006", 5),
    @("003714", "英大睿盛灵活配置混合C", "2.19", "93.65", "7.09", "0.1553", 5),
    @("012202", "中加消费优选混合A", "3.52", "73.15", "3.47", "0.1221", 7),
    @("001678", "英大国企改革主题股票", "1.55", "93.30", "7.23", "0.1121", 2),
    @("004634", "新疆前海联合泳涛灵活配置混合A", "1.20", "92.31", "3.99", "0.0479", 10),
    @("001607", "英大策略优选混合A", "0.57", "91.98", "6.93", "0.0395", 2),
    @("012522", "英大稳固增强核心一年持有混合C", "1.24", "27.71", "2.22", "0.0275", 2),
    @("012203", "中加消费优选混合C", "0.57", "73.15", "3.47", "0.0198", 7),
    @("004436", "汇添富年年泰定期开放混合A", "1.58", "25.75", "1.20", "0.0190", 8),
    @("003447", "英大睿鑫灵活配置混合C", "0.21", "92.71", "8.88", "0.0186", 1),
    @("012521", "英大稳固增强核心一年持有混合A", "0.75", "27.71", "2.22", "0.0166", 2),
    @("007041", "新疆前海联合泳涛灵活配置混合C", "0.31", "92.31", "3.99", "0.0124", 10),
    @("003446", "英大睿鑫灵活配置混合A", "0.07", "92.71", "8.88", "0.0062", 1),
    @("004437", "汇添富年年泰定期开放混合C", "0.14", "25.75", "1.20", "0.0017", 8),
    @("001608", "英大策略优选混合C", "0.02", "91.98", "6.93", "0.0014", 2)
)

function Set-TextValue($cell, $text) {
    # Force a numeric-looking string (fund code "003713", percentages
    # "93.30", etc.) to be written as literal text - preserving leading /
    # trailing zeros - without leaving the cell's number format/style
    # changed afterwards.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

for ($i = 0; $i -lt $q3.Length; $i++) {
    $row = $i + 2
    $rec = $q3[$i]
    $newWs.Cells.Item($row, 1).Value = $i
    Set-TextValue $newWs.Cells.Item($row, 2) $rec[0]
    $newWs.Cells.Item($row, 3).Value = $rec[1]
    Set-TextValue $newWs.Cells.Item($row, 4) $rec[2]
    Set-TextValue $newWs.Cells.Item($row, 5) $rec[3]
    Set-TextValue $newWs.Cells.Item($row, 6) $rec[4]
    Set-TextValue $newWs.Cells.Item($row, 7) $rec[5]
    $newWs.Cells.Item($row, 8).Value = $rec[6]
}

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert the 2022-Q3 row at the top
#    of the data (row 2) and push the rest of the quarters down by one.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$quarters = @(
    @("2022-Q3", 15, 0.8),
    @("2022-Q2", 23, 3.45),
    @("2022-Q1", 14, 1.01),
    @("2021-Q4", 10, 8.83),
    @("2021-Q3", 21, 13.64),
    @("2021-Q2", 30, 18.72),
    @("2021-Q1", 18, 3.67),
    @("2020-Q4", 5, 1.64)
)

# Copy the formatting of the last existing row down onto the brand-new
# row 9 before filling in values, so column A keeps its bold/bordered
# style.
$summary.Cells.Item(8, 1).Copy($summary.Cells.Item(9, 1))

for ($i = 0; $i -lt $quarters.Length; $i++) {
    $row = $i + 2
    $rec = $quarters[$i]
    $summary.Cells.Item($row, 1).Value = $i
    $summary.Cells.Item($row, 2).Value = $rec[0]
    $summary.Cells.Item($row, 3).Value = $rec[1]
    $summary.Cells.Item($row, 4).Value = $rec[2]
}

# Restore the original active sheet/selection so the saved workbook view
# state matches what it was before the edit.
$summary.Activate() | Out-Null
$summary.Range("A1").Select() | Out-Null
